# Add two new weekly timesheet blocks (rows 184-196 and 199-211) to the
# "Spring" sheet, mirroring the existing weekly block pattern (e.g. rows
# 169-181), then adjust the handful of cells whose values differ from the
# template week, and finally re-label one task row with a brand new task
# name ("Project Expo").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spring")

# ---------------------------------------------------------------------
# 1. Clone the most recent weekly block (rows 170-181, the date-header row
#    plus the 11 task/total rows) into the two new weekly slots. This
#    carries over every cell value, number, and style index untouched, so
#    the only remaining work is a handful of per-week value tweaks below.
# ---------------------------------------------------------------------
$ws.Range("A170:I181").Copy($ws.Range("A185"))
$ws.Range("A170:I181").Copy($ws.Range("A200"))

# ---------------------------------------------------------------------
# 2. Re-create the section-divider row ("Date" banner, merged B:H) for
#    each new block. The merge is applied to the still-blank row first,
#    and the look of another divider row (e.g. row 4) is then brought in
#    with a formats-only paste - pasting formats (rather than re-merging
#    a pre-formatted range) keeps every style index identical to the
#    other divider rows instead of minting new border variants for the
#    merged span.
# ---------------------------------------------------------------------
$ws.Range("B184:H184").Merge()
$ws.Range("A4:I4").Copy()
$ws.Range("A184:I184").PasteSpecial(-4122)
$ws.Range("B184").Value = "Date"

$ws.Range("B199:H199").Merge()
$ws.Range("A4:I4").Copy()
$ws.Range("A199:I199").PasteSpecial(-4122)
$ws.Range("B199").Value = "Date"

# ---------------------------------------------------------------------
# 3. Week of 11/27/2023 (rows 185-196) - update the date header and the
#    task rows that differ from the cloned template week.
# ---------------------------------------------------------------------
$ws.Range("B185").Value = 45396
$ws.Range("C185").Value = 45397
$ws.Range("D185").Value = 45398
$ws.Range("E185").Value = 45033
$ws.Range("F185").Value = 45400
$ws.Range("G185").Value = 45401
$ws.Range("H185").Value = 45402

$ws.Range("F193").Value = 2
$ws.Range("G193").Value = 5
$ws.Range("I193").Value = 7

$ws.Range("F194").ClearContents()
$ws.Range("G194").ClearContents()
$ws.Range("H194").ClearContents()
$ws.Range("I194").Value = 5

$ws.Range("F196").Value = 2
$ws.Range("H196").ClearContents()
$ws.Range("I196").Value = 14

# ---------------------------------------------------------------------
# 4. Week of 12/4/2023 (rows 200-211) - same treatment, plus the new
#    "Project Expo" task which replaces the unused "Task zzzzzz" slot.
# ---------------------------------------------------------------------
$ws.Range("B200").Value = 45403
$ws.Range("C200").Value = 45404
$ws.Range("D200").Value = 45405
$ws.Range("E200").Value = 45040
$ws.Range("F200").Value = 45407
$ws.Range("G200").Value = 45408
$ws.Range("H200").Value = 45409

$ws.Range("D208").Value = 2
$ws.Range("E208").Value = 4
$ws.Range("F208").Value = 2
$ws.Range("I208").Value = 6

$ws.Range("E209").ClearContents()
$ws.Range("F209").ClearContents()
$ws.Range("G209").ClearContents()
$ws.Range("H209").ClearContents()
$ws.Range("I209").Value = 1

$ws.Range("A210").Value = "Project Expo"
$ws.Range("G210").Value = 7
$ws.Range("I210").Value = 7

$ws.Range("D211").Value = 2
$ws.Range("F211").Value = 2
$ws.Range("G211").Value = 7
$ws.Range("H211").Value = 0
$ws.Range("I211").Value = 16

# ---------------------------------------------------------------------
# 5. Sheet-level bookkeeping to match the new bottom-of-sheet extent:
#    scroll the view down to the new last block and select its last cell.
# ---------------------------------------------------------------------
$excel.Goto($ws.Range("A201"), $true)
$ws.Range("I211").Select()
